$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the list date (A1)
$ws.Range("A1").Value = 45309

# Update price column (D) for step 1 and step 2 price lists
$ws.Range("D25").Value = 4186.328
$ws.Range("D26").Value = 5159.521
$ws.Range("D27").Value = 5646.136
$ws.Range("D28").Value = 6894.853
$ws.Range("D33").Value = 4723.01
$ws.Range("D34").Value = 5671.219
$ws.Range("D35").Value = 6469.107
$ws.Range("D36").Value = 7478.078
$ws.Range("D37").Value = 8852.048000000001
$ws.Range("D38").Value = 10283.303
$ws.Range("D39").Value = 11879.107
$ws.Range("D40").Value = 14204.868
$ws.Range("D45").Value = 5867.964
$ws.Range("D46").Value = 6941.38
$ws.Range("D47").Value = 8376.178
$ws.Range("D48").Value = 9220.582
$ws.Range("D49").Value = 11360.241
$ws.Range("D50").Value = 12612.587
$ws.Range("D51").Value = 13832.698
$ws.Range("D52").Value = 15779.122
$ws.Range("D53").Value = 17782.853
$ws.Range("D58").Value = 7911.047
$ws.Range("D59").Value = 9220.582
$ws.Range("D60").Value = 10340.505
$ws.Range("D61").Value = 12165.299
$ws.Range("D62").Value = 14025.903
$ws.Range("D63").Value = 16297.947
$ws.Range("D64").Value = 18212.202
$ws.Range("D65").Value = 19428.72
$ws.Range("D66").Value = 22076.494
$ws.Range("D67").Value = 23793.987
$ws.Range("D68").Value = 25547.144
$ws.Range("D69").Value = 28946.302
$ws.Range("D74").Value = 11735.955
$ws.Range("D75").Value = 13131.41
$ws.Range("D76").Value = 14507.831
$ws.Range("D77").Value = 19178.313
$ws.Range("D78").Value = 21826.021
$ws.Range("D79").Value = 23733.114
$ws.Range("D80").Value = 25976.531
$ws.Range("D81").Value = 28345.185
$ws.Range("D82").Value = 30914.23
$ws.Range("D83").Value = 33647.88
$ws.Range("D84").Value = 38964.82
$ws.Range("D85").Value = 44324.653
$ws.Range("D90").Value = 14419.473
$ws.Range("D91").Value = 16154.816
$ws.Range("D92").Value = 19965.433
$ws.Range("D93").Value = 22756.32
$ws.Range("D94").Value = 25547.144
$ws.Range("D95").Value = 28409.612
$ws.Range("D96").Value = 31236.24
$ws.Range("D97").Value = 34456.444
$ws.Range("D98").Value = 38034.502
$ws.Range("D99").Value = 40789.588
$ws.Range("D100").Value = 46908.057
$ws.Range("D101").Value = 53062.258
$ws.Range("D102").Value = 59252.254
$ws.Range("D103").Value = 65442.232
$ws.Range("D108").Value = 34993.184
$ws.Range("D109").Value = 39251.082
$ws.Range("D110").Value = 43437.327
$ws.Range("D111").Value = 47945.629
$ws.Range("D112").Value = 52357.413
$ws.Range("D113").Value = 61220.166
$ws.Range("D114").Value = 70129.503
$ws.Range("D115").Value = 79074.58
$ws.Range("D116").Value = 88413.25199999999
$ws.Range("D117").Value = 98145.47199999999
$ws.Range("D118").Value = 118075.151
$ws.Range("D119").Value = 129167.07
$ws.Range("D120").Value = 140258.972
$ws.Range("D121").Value = 151708.67
$ws.Range("D122").Value = 177076.943
